# Applies corrections described by the commit message
# "Correcion a Diebold Mariano y revision de Cap1"
#
# Updates column K (EnCQR-LSTM) values on sheet "Detalle_Pasos" for rows 2-25,
# and column C (Empirical) values on sheet "Reliability_Data" for select rows
# in the 799-892 range.

$wb = $excel.ActiveWorkbook

$wsDetalle = $wb.Worksheets.Item("Detalle_Pasos")
$wsReliability = $wb.Worksheets.Item("Reliability_Data")

# --- Sheet "Detalle_Pasos": column K (EnCQR-LSTM), rows 2-25 ---
$kValues = @{
    2  = 3.758027037470606
    3  = 3.423949234137423
    4  = 3.749313837658861
    5  = 3.79187289159037
    6  = 3.796472818008077
    7  = 3.88462698097165
    8  = 3.875022854871307
    9  = 3.842861328867081
    10 = 3.838283094437498
    11 = 3.907072762435249
    12 = 3.865202952844662
    13 = 3.621952348279919
    14 = 4.364350586382203
    15 = 11.18854862243428
    16 = 12.93301175572414
    17 = 4.230612849622545
    18 = 3.593957573726737
    19 = 3.889728855548245
    20 = 3.403864253255919
    21 = 2.96815869810564
    22 = 3.198703415617933
    23 = 3.407760641777841
    24 = 3.456431200976954
    25 = 3.557904840160039
}

foreach ($row in $kValues.Keys) {
    $wsDetalle.Range("K$row").Value = $kValues[$row]
}

# --- Sheet "Reliability_Data": column C (Empirical), selected rows ---
$cValues = @{
    799 = 0.04166666666666666
    800 = 0.04166666666666666
    801 = 0.04166666666666666
    802 = 0.04166666666666666
    803 = 0.04166666666666666
    804 = 0.04166666666666666
    805 = 0.04166666666666666
    806 = 0.04166666666666666
    807 = 0.04166666666666666
    808 = 0.04166666666666666
    824 = 0.08333333333333333
    825 = 0.08333333333333333
    826 = 0.08333333333333333
    827 = 0.08333333333333333
    828 = 0.08333333333333333
    829 = 0.125
    830 = 0.1666666666666667
    831 = 0.1666666666666667
    832 = 0.2083333333333333
    833 = 0.2083333333333333
    834 = 0.2083333333333333
    835 = 0.2083333333333333
    841 = 0.3333333333333333
    842 = 0.4166666666666667
    843 = 0.4583333333333333
    844 = 0.4583333333333333
    846 = 0.5
    847 = 0.7083333333333334
    848 = 0.9166666666666666
    852 = 0.9583333333333334
    853 = 0.9583333333333334
    854 = 0.9583333333333334
    855 = 0.9583333333333334
    856 = 0.9583333333333334
    857 = 0.9583333333333334
    858 = 0.9583333333333334
    859 = 0.9583333333333334
    860 = 0.9583333333333334
    889 = 1
    890 = 1
    891 = 1
    892 = 1
}

foreach ($row in $cValues.Keys) {
    $wsReliability.Range("C$row").Value = $cValues[$row]
}
